# Update cryptos list data per the scraped diff (commit: Tue Oct 22 23:50:23 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. '1.00', '593.72') need the
# cell's number format forced to Text ('@') first, otherwise Excel's COM layer
# auto-converts the assigned string into a numeric value (losing the original
# formatted-as-text representation, e.g. '1.00' -> 1).

$ws.Range("D2").Value = '67.358.04'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '2.620.17'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.72'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.21'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("D9").Value = '2.619.76'
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("E10").Value = '  -2.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.364'
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.22'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.63'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '3.117.57'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000181'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").Value = '67.261.65'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '2.628.62'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.99'
$ws.Range("E19").Value = '  +2.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.01'
$ws.Range("E20").Value = '  +2.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '357.48'
$ws.Range("E21").Value = '  -1.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.68'
$ws.Range("E23").Value = '  -2.06%  '
$ws.Range("E24").Value = '  -3.96%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.35'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '69.92'
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("E28").Value = '  -1.72%  '
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '545.44'
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.95'
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.35'
$ws.Range("E33").Value = '  -2.52%  '
$ws.Range("E34").Value = '  -1.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.136'
$ws.Range("E35").Value = '  +4.77%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.81'
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.04'
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.366'
$ws.Range("E40").Value = '  -1.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.27'
$ws.Range("E41").Value = '  +1.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.20'
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").Value = '  -3.30%  '
$ws.Range("D46").Value = '0.0₆0301'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '152.57'
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.580'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.80'
$ws.Range("E49").Value = '  -1.52%  '
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("E51").Value = '  -0.49%  '
